# Refresh the crypto price/volume snapshot (columns D = Price, E = Volume(1h))
# for the rows whose source data changed, matching the latest scrape.
# Price cells that look numeric ("1.002", "0.5227", ...) are written with a
# leading apostrophe so Excel keeps them as text (same as the original cells)
# instead of silently re-typing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, new Price (D) value [or $null = unchanged], new Volume(1h) (E) value [or $null = unchanged]
$updates = @(
    [PSCustomObject]@{ Row = 2; Price = "27.124.88"; Volume = "  -0.36%  " }
    [PSCustomObject]@{ Row = 3; Price = "1.891.41"; Volume = "  -0.84%  " }
    [PSCustomObject]@{ Row = 4; Price = $null; Volume = "  +0.18%  " }
    [PSCustomObject]@{ Row = 5; Price = "306.88"; Volume = "  -0.36%  " }
    [PSCustomObject]@{ Row = 6; Price = "1.002"; Volume = "  +0.13%  " }
    [PSCustomObject]@{ Row = 7; Price = "0.5227"; Volume = "  -0.48%  " }
    [PSCustomObject]@{ Row = 8; Price = "0.3756"; Volume = "  -0.69%  " }
    [PSCustomObject]@{ Row = 9; Price = $null; Volume = "  -0.23%  " }
    [PSCustomObject]@{ Row = 10; Price = "21.07"; Volume = "  -0.85%  " }
    [PSCustomObject]@{ Row = 11; Price = "0.8993"; Volume = "  +0.18%  " }
    [PSCustomObject]@{ Row = 12; Price = $null; Volume = "  +6.10%  " }
    [PSCustomObject]@{ Row = 13; Price = "1.934.60"; Volume = "  +1.55%  " }
    [PSCustomObject]@{ Row = 14; Price = "96.03"; Volume = "  +0.89%  " }
    [PSCustomObject]@{ Row = 15; Price = "5.286"; Volume = "  +0.19%  " }
    [PSCustomObject]@{ Row = 16; Price = "1.003"; Volume = "  +0.16%  " }
    [PSCustomObject]@{ Row = 17; Price = $null; Volume = "  -0.77%  " }
    [PSCustomObject]@{ Row = 18; Price = $null; Volume = "  +0.58%  " }
    [PSCustomObject]@{ Row = 19; Price = $null; Volume = "  +0.20%  " }
    [PSCustomObject]@{ Row = 20; Price = "27.149.92"; Volume = "  -0.48%  " }
    [PSCustomObject]@{ Row = 21; Price = "5.081"; Volume = "  +0.08%  " }
    [PSCustomObject]@{ Row = 22; Price = "10.68"; Volume = "  +0.46%  " }
    [PSCustomObject]@{ Row = 23; Price = "6.400"; Volume = "  -0.65%  " }
    [PSCustomObject]@{ Row = 24; Price = "147.91"; Volume = "  +1.35%  " }
    [PSCustomObject]@{ Row = 25; Price = "2.282"; Volume = "  -2.03%  " }
    [PSCustomObject]@{ Row = 26; Price = $null; Volume = "  +0.17%  " }
    [PSCustomObject]@{ Row = 27; Price = "1.740"; Volume = "  +0.00%  " }
    [PSCustomObject]@{ Row = 28; Price = "114.98"; Volume = "  +0.12%  " }
    [PSCustomObject]@{ Row = 29; Price = "4.784"; Volume = "  -0.72%  " }
    [PSCustomObject]@{ Row = 30; Price = "4.846"; Volume = "  -2.44%  " }
    [PSCustomObject]@{ Row = 31; Price = "0.09227"; Volume = "  -0.05%  " }
    [PSCustomObject]@{ Row = 32; Price = "0.05035"; Volume = "  -0.71%  " }
    [PSCustomObject]@{ Row = 33; Price = "0.7872"; Volume = "  -4.13%  " }
    [PSCustomObject]@{ Row = 34; Price = "1.212"; Volume = "  -2.63%  " }
    [PSCustomObject]@{ Row = 35; Price = "3.423"; Volume = "  +3.61%  " }
    [PSCustomObject]@{ Row = 36; Price = "2.957"; Volume = "  -1.08%  " }
    [PSCustomObject]@{ Row = 37; Price = "2.584"; Volume = "  -0.64%  " }
    [PSCustomObject]@{ Row = 38; Price = "0.5685"; Volume = "  +0.19%  " }
    [PSCustomObject]@{ Row = 39; Price = "0.01977"; Volume = "  -0.60%  " }
    [PSCustomObject]@{ Row = 40; Price = "1.075"; Volume = "  -0.01%  " }
    [PSCustomObject]@{ Row = 41; Price = "9.013"; Volume = "  +0.21%  " }
    [PSCustomObject]@{ Row = 42; Price = "6.550"; Volume = "  -1.19%  " }
    [PSCustomObject]@{ Row = 43; Price = "116.26"; Volume = "  -2.61%  " }
    [PSCustomObject]@{ Row = 44; Price = "0.1514"; Volume = "  -0.13%  " }
    [PSCustomObject]@{ Row = 45; Price = "0.4851"; Volume = $null }
    [PSCustomObject]@{ Row = 46; Price = "1.002"; Volume = "  +0.10%  " }
    [PSCustomObject]@{ Row = 47; Price = "10.07"; Volume = "  -1.13%  " }
    [PSCustomObject]@{ Row = 48; Price = "1.618"; Volume = "  -0.56%  " }
    [PSCustomObject]@{ Row = 49; Price = $null; Volume = "  +1.30%  " }
    [PSCustomObject]@{ Row = 50; Price = "63.47"; Volume = "  -0.44%  " }
    [PSCustomObject]@{ Row = 51; Price = "0.05937"; Volume = "  +0.05%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Cells.Item($u.Row, 4)   # column D
        $priceValue = $u.Price
        $isNumericLooking = $priceValue -match '^[+-]?\d+(\.\d+)?$'
        if ($isNumericLooking) {
            # Force text storage (matches original inlineStr cell) instead of
            # letting Excel auto-convert the digits to a real number.
            $priceCell.Value = "'" + $priceValue
        } else {
            $priceCell.Value = $priceValue
        }
    }
    if ($null -ne $u.Volume) {
        $ws.Cells.Item($u.Row, 5).Value = $u.Volume   # column E
    }
}

